$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.175.41"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "2.251.56"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'307.83"
$ws.Range("E5").Value = "  -4.59%  "
$ws.Range("D6").Value = "'98.37"
$ws.Range("E6").Value = "  -2.98%  "
$ws.Range("E7").Value = "  -0.70%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("E9").Value = "  -4.09%  "
$ws.Range("D10").Value = "'35.52"
$ws.Range("E10").Value = "  -3.92%  "
$ws.Range("D11").Value = "'0.0820"
$ws.Range("E11").Value = "  -1.36%  "
$ws.Range("D12").Value = "'7.30"
$ws.Range("E12").Value = "  -5.72%  "
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("D14").Value = "2.594.80"
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("D15").Value = "2.253.93"
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("D16").Value = "'0.835"
$ws.Range("E16").Value = "  -2.38%  "
$ws.Range("D17").Value = "'13.77"
$ws.Range("E17").Value = "  -2.48%  "
$ws.Range("D18").Value = "44.005.43"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("D19").Value = "'12.78"
$ws.Range("E19").Value = "  -6.03%  "
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("E21").Value = "  -3.94%  "
$ws.Range("D22").Value = "'65.30"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "'241.63"
$ws.Range("E23").Value = "  +1.89%  "
$ws.Range("E24").Value = "  -7.41%  "
$ws.Range("D25").Value = "'1.97"
$ws.Range("E25").Value = "  -8.52%  "
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").Value = "  -2.09%  "
$ws.Range("D29").Value = "'36.60"
$ws.Range("E29").Value = "  -0.63%  "
$ws.Range("D30").Value = "'6.19"
$ws.Range("E30").Value = "  -1.82%  "
$ws.Range("D31").Value = "'20.13"
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("D32").Value = "'156.88"
$ws.Range("E32").Value = "  -2.30%  "
$ws.Range("E33").Value = "  +14.15%  "
$ws.Range("D34").Value = "'0.0819"
$ws.Range("E34").Value = "  -4.10%  "
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("E37").Value = "  -4.68%  "
$ws.Range("E38").Value = "  -3.69%  "
$ws.Range("D39").Value = "'15.46"
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("D40").Value = "'3.85"
$ws.Range("E40").Value = "  -9.18%  "
$ws.Range("D41").Value = "'0.0306"
$ws.Range("E41").Value = "  -3.67%  "
$ws.Range("E42").Value = "  -10.59%  "
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").Value = "1.760.93"
$ws.Range("E44").Value = "  -2.00%  "
$ws.Range("D45").Value = "'86.62"
$ws.Range("E45").Value = "  +5.11%  "
$ws.Range("E46").Value = "  -0.72%  "
$ws.Range("D47").Value = "'0.192"
$ws.Range("E47").Value = "  -3.40%  "
$ws.Range("D48").Value = "'101.16"
$ws.Range("E48").Value = "  -2.35%  "
$ws.Range("D49").Value = "'8.24"
$ws.Range("E49").Value = "  -2.20%  "
$ws.Range("D50").Value = "'55.54"
$ws.Range("E50").Value = "  -5.49%  "
$ws.Range("D51").Value = "'69.45"
$ws.Range("E51").Value = "  -7.66%  "
